$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 26 for the newly added stock (105560 / KB금융).
# This shifts rows 26-29 down to 27-30, carrying their formatting (styles) along.
$ws.Rows(26).Insert()

# --- New row 26: 105560 / KB금융 ---
$ws.Range("A26").Value = "'105560"
$ws.Range("B26").Value = "KB금융"
$ws.Range("C26").Value = 111500
$ws.Range("D26").Value = 0.0362
$ws.Range("E26").Value = 8.86
$ws.Range("F26").Value = 11.69
$ws.Range("G26").Value = 8.51
$ws.Range("H26").Value = 152301.9
$ws.Range("I26").Value = 2.85
$ws.Range("J26").Value = 96
$ws.Range("K26").Value = 96

# --- Refresh market-data columns (price, change%, dividend yield, stochastic K/D) ---
# rows 2-25 keep their original row position; rows that were 26-29 are now 27-30.
$ws.Range("C2").Value = 19270
$ws.Range("D2").Value = 0.018
$ws.Range("I2").Value = 5.19
$ws.Range("J2").Value = 67
$ws.Range("K2").Value = 67

$ws.Range("C3").Value = 97100
$ws.Range("D3").Value = 0.0125
$ws.Range("I3").Value = 6.69
$ws.Range("J3").Value = 63
$ws.Range("K3").Value = 63

$ws.Range("C4").Value = 438000
$ws.Range("D4").Value = 0.0294
$ws.Range("I4").Value = 4.34
$ws.Range("J4").Value = 91
$ws.Range("K4").Value = 91

$ws.Range("C5").Value = 31500
$ws.Range("D5").Value = 0.0096
$ws.Range("I5").Value = 6.35
$ws.Range("J5").Value = 49
$ws.Range("K5").Value = 49

$ws.Range("C6").Value = 26900
$ws.Range("D6").Value = 0.076
$ws.Range("I6").Value = 4.46
$ws.Range("J6").Value = 65
$ws.Range("K6").Value = 65

$ws.Range("C7").Value = 24950
$ws.Range("D7").Value = 0.0163
$ws.Range("I7").Value = 4.81
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 90

$ws.Range("C8").Value = 10240
$ws.Range("D8").Value = -0.0039
$ws.Range("I8").Value = 5.03
$ws.Range("J8").Value = 99
$ws.Range("K8").Value = 99

$ws.Range("C9").Value = 71400
$ws.Range("D9").Value = 0.0469
$ws.Range("I9").Value = 4.2
$ws.Range("J9").Value = 54
$ws.Range("K9").Value = 54

$ws.Range("C10").Value = 205500
$ws.Range("D10").Value = 0.0199
$ws.Range("I10").Value = 5.84
$ws.Range("J10").Value = 53
$ws.Range("K10").Value = 53

$ws.Range("C11").Value = 121800
$ws.Range("D11").Value = 0.0227
$ws.Range("I11").Value = 5.58
$ws.Range("J11").Value = 97
$ws.Range("K11").Value = 97

$ws.Range("C12").Value = 20650
$ws.Range("D12").Value = 0.0633
$ws.Range("I12").Value = 4.6
$ws.Range("J12").Value = 97
$ws.Range("K12").Value = 97

$ws.Range("C13").Value = 75000
$ws.Range("D13").Value = 0.049
$ws.Range("I13").Value = 4.67
$ws.Range("J13").Value = 97
$ws.Range("K13").Value = 97

$ws.Range("C14").Value = 57000
$ws.Range("D14").Value = 0.016
$ws.Range("I14").Value = 6.21
$ws.Range("J14").Value = 78
$ws.Range("K14").Value = 78

$ws.Range("C15").Value = 85700
$ws.Range("D15").Value = 0.0263
$ws.Range("I15").Value = 6.42
$ws.Range("J15").Value = 91
$ws.Range("K15").Value = 91

$ws.Range("C16").Value = 17960
$ws.Range("D16").Value = 0.031
$ws.Range("I16").Value = 5.93
$ws.Range("J16").Value = 98
$ws.Range("K16").Value = 98

$ws.Range("C17").Value = 49150
$ws.Range("D17").Value = 0.002
$ws.Range("I17").Value = 5.7
$ws.Range("J17").Value = 91
$ws.Range("K17").Value = 91

$ws.Range("C18").Value = 20450
$ws.Range("D18").NumberFormat = "0.00%"
$ws.Range("D18").Value = 0.0049
$ws.Range("I18").Value = 6.01
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 40

$ws.Range("C19").Value = 52100
$ws.Range("D19").Value = 0.0216
$ws.Range("I19").Value = 3.84
$ws.Range("J19").Value = 93
$ws.Range("K19").Value = 93

$ws.Range("C20").Value = 13750
$ws.Range("D20").Value = 0.0231
$ws.Range("I20").Value = 4.73
$ws.Range("J20").Value = 63
$ws.Range("K20").Value = 63

$ws.Range("C21").Value = 124100
$ws.Range("D21").Value = 0.0032
$ws.Range("I21").Value = 4.35
$ws.Range("J21").Value = 93
$ws.Range("K21").Value = 93

$ws.Range("C22").Value = 38900
$ws.Range("D22").Value = 0.0291
$ws.Range("I22").Value = 3.74
$ws.Range("J22").Value = 39
$ws.Range("K22").Value = 39

$ws.Range("C23").Value = 61300
$ws.Range("D23").Value = 0.02
$ws.Range("I23").Value = 3.52
$ws.Range("J23").Value = 91
$ws.Range("K23").Value = 91

$ws.Range("C24").Value = 46450
$ws.Range("D24").Value = 0.0087
$ws.Range("I24").Value = 5.81
$ws.Range("J24").Value = 59
$ws.Range("K24").Value = 59

$ws.Range("C25").Value = 81800
$ws.Range("D25").Value = 0.0251
$ws.Range("I25").Value = 4.4
$ws.Range("J25").Value = 98
$ws.Range("K25").Value = 98

$ws.Range("C27").Value = 12110
$ws.Range("D27").Value = 0.0058
$ws.Range("I27").Value = 5.37
$ws.Range("J27").Value = 95
$ws.Range("K27").Value = 95

$ws.Range("C28").Value = 12090
$ws.Range("D28").Value = 0.016
$ws.Range("I28").Value = 4.14
$ws.Range("J28").Value = 98
$ws.Range("K28").Value = 98

$ws.Range("C29").Value = 21550
$ws.Range("D29").NumberFormat = "0.00%"
$ws.Range("D29").Value = 0.0262
$ws.Range("I29").Value = 4.62
$ws.Range("J29").Value = 96
$ws.Range("K29").Value = 96

$ws.Range("C30").Value = 22100
$ws.Range("D30").Value = 0.0351
$ws.Range("I30").Value = 5.43
$ws.Range("J30").Value = 99
$ws.Range("K30").Value = 99

# --- View state: scroll so row 8 is at top, select I24 (as in the saved file) ---
$ws.Range("I24").Select()
$excel.ActiveWindow.ScrollRow = 8
